# Insert a new data row at row 4 (pushing existing rows 4..74 down to 5..75)
# and populate it with the new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(4).Insert()

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C4").Value = "Arica y Parinacota"
$ws.Range("D4").Value = (Get-Date -Year 2022 -Month 7 -Day 7 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E4").Value = 15
$ws.Range("F4").Value = 100114001
$ws.Range("G4").Value = "Papa"
$ws.Range("H4").Value = "Asterix"
$ws.Range("I4").Value = "1a (guarda)"
$ws.Range("J4").Value = 1000
$ws.Range("K4").Value = 10000
$ws.Range("L4").Value = 11000
$ws.Range("M4").Value = 10500
$ws.Range("N4").Value = "`$/saco 25 kilos"
$ws.Range("O4").Value = "Región de Los Lagos"
$ws.Range("P4").Value = 420
$ws.Range("Q4").Value = 25
$ws.Range("R4").Value = "Hortaliza"
